# Auto-generated edit script: updates numeric cells in the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to refresh computed market-board
# price/profit figures per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4879.5835
$ws.Range("J32").Value = 5105.091
$ws.Range("L32").Value = 5105.091
$ws.Range("N32").Value = -5757.091
$ws.Range("H86").Value = 1091.6086
$ws.Range("I86").Value = 835.58826
$ws.Range("J86").Value = 1817
$ws.Range("K86").Value = 835.58826
$ws.Range("L86").Value = 1817
$ws.Range("M86").Value = 287.41174
$ws.Range("N86").Value = -4063
$ws.Range("H89").Value = 1091.6086
$ws.Range("I89").Value = 835.58826
$ws.Range("J89").Value = 1817
$ws.Range("K89").Value = 4177.9413
$ws.Range("L89").Value = 9085
$ws.Range("M89").Value = 1438.0587
$ws.Range("N89").Value = -20317
$ws.Range("H98").Value = 4242.0586
$ws.Range("I98").Value = 3861.4546
$ws.Range("J98").Value = 4939.8335
$ws.Range("K98").Value = 3861.4546
$ws.Range("L98").Value = 4939.8335
$ws.Range("M98").Value = -2363.4546
$ws.Range("N98").Value = -7935.8335
$ws.Range("H107").Value = 1299.6666
$ws.Range("J107").Value = 400
$ws.Range("L107").Value = 400
$ws.Range("N107").Value = -4240
$ws.Range("H122").Value = 4242.0586
$ws.Range("I122").Value = 3861.4546
$ws.Range("J122").Value = 4939.8335
$ws.Range("K122").Value = 11584.3638
$ws.Range("L122").Value = 14819.5005
$ws.Range("M122").Value = -9134.363799999999
$ws.Range("N122").Value = -19719.5005
$ws.Range("H139").Value = 149499.75
$ws.Range("J139").Value = 159333
$ws.Range("L139").Value = 159333
$ws.Range("N139").Value = -169613
$ws.Range("H140").Value = 139998.5
$ws.Range("J140").Value = 149997
$ws.Range("L140").Value = 149997
$ws.Range("N140").Value = -160357

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4090.6516
$ws.Range("I32").Value = 3568.224
$ws.Range("J32").Value = 7878.25
$ws.Range("K32").Value = 3568.224
$ws.Range("L32").Value = 7878.25
$ws.Range("M32").Value = -3281.224
$ws.Range("N32").Value = -8452.25
$ws.Range("H80").Value = 76663.336
$ws.Range("J80").Value = 76663.336
$ws.Range("L80").Value = 76663.336
$ws.Range("N80").Value = -78659.336
$ws.Range("H83").Value = 76663.336
$ws.Range("J83").Value = 76663.336
$ws.Range("L83").Value = 229990.008
$ws.Range("N83").Value = -239974.008
$ws.Range("H109").Value = 25500
$ws.Range("J109").Value = 25500
$ws.Range("L109").Value = 25500
$ws.Range("N109").Value = -28274
$ws.Range("H122").Value = 20835888
$ws.Range("I122").Value = 2420.6667
$ws.Range("K122").Value = 7262.000100000001
$ws.Range("M122").Value = -4812.000100000001
$ws.Range("H134").Value = 85069
$ws.Range("J134").Value = 85069
$ws.Range("L134").Value = 85069
$ws.Range("N134").Value = -95209
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6553.5947
$ws.Range("I86").Value = 7661.4346
$ws.Range("J86").Value = 4733.5713
$ws.Range("K86").Value = 7661.4346
$ws.Range("L86").Value = 4733.5713
$ws.Range("M86").Value = -6538.4346
$ws.Range("N86").Value = -6979.5713
$ws.Range("H89").Value = 6553.5947
$ws.Range("I89").Value = 7661.4346
$ws.Range("J89").Value = 4733.5713
$ws.Range("K89").Value = 38307.173
$ws.Range("L89").Value = 23667.8565
$ws.Range("M89").Value = -32691.173
$ws.Range("N89").Value = -34899.85649999999
$ws.Range("H99").Value = 4210
$ws.Range("I99").Value = 2703.5
$ws.Range("K99").Value = 2703.5
$ws.Range("M99").Value = -1205.5
$ws.Range("H107").Value = 2657.3333
$ws.Range("I107").Value = 2657.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2657.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -737.3332999999998
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 4702.9414
$ws.Range("I134").Value = 4662.2856
$ws.Range("K134").Value = 13986.8568
$ws.Range("M134").Value = -11451.8568
$ws.Range("H137").Value = 200000
$ws.Range("J137").Value = 200000
$ws.Range("L137").Value = 200000
$ws.Range("N137").Value = -210200
$ws.Range("H138").Value = 72774
$ws.Range("I138").Value = 78999
$ws.Range("J138").Value = 69661.5
$ws.Range("K138").Value = 78999
$ws.Range("L138").Value = 69661.5
$ws.Range("M138").Value = -73859
$ws.Range("N138").Value = -79941.5
$ws.Range("H140").Value = 189998.5
$ws.Range("J140").Value = 189998.5
$ws.Range("L140").Value = 189998.5
$ws.Range("N140").Value = -200358.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2506.1667
$ws.Range("I134").Value = 1793.9286
$ws.Range("K134").Value = 5381.7858
$ws.Range("M134").Value = -2846.7858
$ws.Range("H138").Value = 96185.09
$ws.Range("J138").Value = 73115.11
$ws.Range("L138").Value = 73115.11
$ws.Range("N138").Value = -83395.11
$ws.Range("H140").Value = 85144.52
$ws.Range("J140").Value = 87366.3
$ws.Range("L140").Value = 87366.3
$ws.Range("N140").Value = -97726.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 21874.46
$ws.Range("I56").Value = 21874.46
$ws.Range("K56").Value = 21874.46
$ws.Range("M56").Value = -21344.46
$ws.Range("H122").Value = 1873.2307
$ws.Range("I122").Value = 531.6
$ws.Range("J122").Value = 2711.75
$ws.Range("K122").Value = 4784.400000000001
$ws.Range("L122").Value = 24405.75
$ws.Range("M122").Value = -2334.400000000001
$ws.Range("N122").Value = -29305.75
$ws.Range("H131").Value = 1268.2
$ws.Range("I131").Value = 898.1539
$ws.Range("J131").Value = 1669.0834
$ws.Range("K131").Value = 2694.4617
$ws.Range("L131").Value = 5007.2502
$ws.Range("M131").Value = 2345.5383
$ws.Range("N131").Value = -15087.2502

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2778079
$ws.Range("I2").Value = 7142933
$ws.Range("K2").Value = 7142933
$ws.Range("M2").Value = -7142820
$ws.Range("H122").Value = 55557320
$ws.Range("I122").Value = 2279.6
$ws.Range("J122").Value = 125001130
$ws.Range("K122").Value = 6838.799999999999
$ws.Range("L122").Value = 375003390
$ws.Range("M122").Value = -4388.799999999999
$ws.Range("N122").Value = -375008290
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1857.421
$ws.Range("I22").Value = 1880.2
$ws.Range("K22").Value = 1880.2
$ws.Range("M22").Value = -1585.2
$ws.Range("H27").Value = 1857.421
$ws.Range("I27").Value = 1880.2
$ws.Range("K27").Value = 1880.2
$ws.Range("M27").Value = -1773.2
$ws.Range("H39").Value = 15750
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -1040
$ws.Range("N39").Value = -30920
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 100000
$ws.Range("I17").Value = 100000
$ws.Range("K17").Value = 100000
$ws.Range("M17").Value = -99828
$ws.Range("H44").Value = 33747.5
$ws.Range("I44").Value = 30000
$ws.Range("K44").Value = 30000
$ws.Range("M44").Value = -29446
$ws.Range("H46").Value = 121666.664
$ws.Range("J46").Value = 121666.664
$ws.Range("L46").Value = 121666.664
$ws.Range("N46").Value = -122128.664
$ws.Range("H58").Value = 30493.8
$ws.Range("I58").Value = 28156.666
$ws.Range("K58").Value = 28156.666
$ws.Range("M58").Value = -27848.666
$ws.Range("H107").Value = 1998.2222
$ws.Range("I107").Value = 1295
$ws.Range("K107").Value = 3885
$ws.Range("M107").Value = -1965
$ws.Range("H134").Value = 121666.664
$ws.Range("J134").Value = 121666.664
$ws.Range("L134").Value = 364999.992
$ws.Range("N134").Value = -370069.992
$ws.Range("H136").Value = 2695.36
$ws.Range("I136").Value = 2019.2
$ws.Range("K136").Value = 6057.6
$ws.Range("M136").Value = -3507.6
$ws.Range("H138").Value = 116666.664
$ws.Range("I138").Value = 100000
$ws.Range("J138").Value = 150000
$ws.Range("K138").Value = 100000
$ws.Range("L138").Value = 150000
$ws.Range("M138").Value = -94860
$ws.Range("N138").Value = -160280

